$d = $word.ActiveDocument
$lb = [char]11

# ---------------------------------------------------------------------------
# The edit reshuffles the "Objetivos" / "Docente(s) Responsavel(eis)" /
# "Programa resumido" / "Programa" / "Avaliacao" / "Bibliografia" sections:
# every paragraph keeps its original style in place, only the *text*
# carried by eight paragraphs gets swapped around (a rotation of content
# between slots). No paragraphs are added or removed and no styles change.
# ---------------------------------------------------------------------------

# Paragraph 6 (plain): was the PT "Familiarizar..." objective text,
# becomes the PT "Programa resumido" summary text.
$d.Paragraphs.Item(6).Range.Text = "Integração de funções de uma variável real. Funções reais de variáveis reais, Diferenciabilidade, Derivada direcional. Máximos e Mínios em domínios abertos e Multiplicadores de Lagrange"

# Paragraph 7 (italic): was the EN "The discipline aims..." objective text,
# becomes the EN "Programa resumido" summary text.
$d.Paragraphs.Item(7).Range.Text = "Integration of real functions. Real functions with several variables, Differentiability, Directional derivatives. Maximum and minimum in open domains, Lagrange Multipliers."

# Paragraph 9 (ListBullet): was "8822123 - Roberta Veloso Garcia",
# becomes the PT "Familiarizar..." objective text.
$d.Paragraphs.Item(9).Range.Text = "Familiarizar os alunos com resultados fundamentais relativos a: integração de funções de uma variável real, cálculo diferencial de funções de n variáveis reais  e suas aplicações."

# Paragraph 11 (plain, multi-line): was the short PT "Programa resumido" text,
# becomes the long PT "Programa" content (six lines joined by manual breaks).
$d.Paragraphs.Item(11).Range.Text = "Integração de funções reais: Primitivas (Integral indefinida), Integral de Riemann (Integral definida), Teorema fundamental do cálculo, Técnicas de integração e aplicações. " + $lb + "O espaço euclidiano R^n: Conjuntos abertos, fechados e compactos." + $lb + "Funções de n várias variáveis Reais: Gráficos e curvas de nível de funções de duas variáveis." + $lb + "Limites e Continuidade: Teorema de Weierstrass" + $lb + "Diferenciabilidade: Derivadas parciais, diferencial total, derivadas parciais de ordem superior, teorema de Schwarz, regra da cadeia, planos tangentes e aproximações lineares, derivada direcional, vetor gradiente, teorema da função implícita, jacobiano." + $lb + "Máximos e mínimos: Valores Extremos de funções de duas ou mais variáveis em domínios abertos, Hessiano de uma função real de n variáveis, multiplicadores de Lagrange."

# Paragraph 12 (italic): was the short EN "Programa resumido" text,
# becomes the EN "The discipline aims..." objective text.
$d.Paragraphs.Item(12).Range.Text = "The discipline aims at familiarizing students with fundamental results regarding: integration of real functions, Differential calculus for functions of n real variables and applications"

# Paragraph 14 (plain, multi-line): was the long PT "Programa" content,
# becomes the "Metodo" evaluation text (single line, no breaks).
$d.Paragraphs.Item(14).Range.Text = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Paragraph 17 ("Avaliação" bullet list: Método / Critério / Norma de
# recuperação). Re-point each value run to the next one's old content,
# shifting Método -> Critério -> Norma -> Bibliografia list. Process from
# the last value to the first so that a freshly written value can never be
# re-matched by a later Find on stale text. A fresh Range is fetched for
# every call because Find.Execute collapses its Range to the (replaced)
# match, so a stale Range would no longer cover the rest of the paragraph.

$d.Paragraphs.Item(17).Range.Find.Execute(
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ("GUIDORIZZI, Hamilton L. UM CURSO DE CÁLCULO,2011, 5. ed., v.2" + $lb +
     "LEITHOLD, Louis. CÁLCULO COM GEOMETRIA ANALÍTICA, São Paulo: HARBRA LTDA, 1990. v.2" + $lb +
     "ANTON, Howard; BIVENS, Irl, DAVIS, Stephen. CÁLCULO, 8. ed. São Paulo:Pearson, 2011, v.2" + $lb +
     "SIMMONS, George F. CÁLCULO COM GEOMETRIA ANALÍTICA, São Paulo: Pearson, 2014. v.2" + $lb +
     "STEWART, James. CÁLCULO. revisão técnica Ricardo Miranda Martins. 7. ed. São Paulo: Cengage Learning, 2013. v.2" + $lb +
     "THOMAS, George B. WEIR, Maurice D.; HASS, Joel; GIORDANO, CÁLCULO. revisão técnica Cláudio Hirofume Asano .12.ed. São Paulo: Pearson Education do Brasil, 2013. v.2"),
    2) | Out-Null

$d.Paragraphs.Item(17).Range.Find.Execute(
    "NF≥ 5,0.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.",
    2) | Out-Null

$d.Paragraphs.Item(17).Range.Find.Execute(
    "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NF≥ 5,0.",
    2) | Out-Null

# Paragraph 19 (plain, multi-line): was the Bibliography list,
# becomes "8822123 - Roberta Veloso Garcia".
$d.Paragraphs.Item(19).Range.Text = "8822123 - Roberta Veloso Garcia"
